$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns E, F, G
# Shared-string insertion order matters: cxq (4), hyy (5), hzj (6)
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

# Data row 2
$ws.Range("E2").Value = 0.91691394658753711
$ws.Range("F2").Value = 0.89789789789789787
$ws.Range("G2").Value = 0.90909090909090906

# Data row 3
$ws.Range("E3").Value = 0.89795918367346939
$ws.Range("F3").Value = 0.88064516129032255
$ws.Range("G3").Value = 0.88395904436860073

# Update selection to match target state (F1 active, full-column selection)
$ws.Range("F1:F1048576").Select()
